$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 404560.12
$ws.Range("I2").Value = 559614.1
$ws.Range("J2").Value = 1419.6
$ws.Range("K2").Value = 559614.1
$ws.Range("L2").Value = 1419.6
$ws.Range("M2").Value = -559501.1
$ws.Range("N2").Value = -1645.6

$ws.Range("H18").Value = 619.8
$ws.Range("I18").Value = 619.8
$ws.Range("K18").Value = 619.8
$ws.Range("M18").Value = -335.8

$ws.Range("H41").Value = 926.2857
$ws.Range("I41").Value = 984.75
$ws.Range("J41").Value = 848.3333
$ws.Range("K41").Value = 984.75
$ws.Range("L41").Value = 848.3333
$ws.Range("M41").Value = -544.75
$ws.Range("N41").Value = -1728.3333

$ws.Range("H64").Value = 5196.8887
$ws.Range("J64").Value = 5326.2856
$ws.Range("L64").Value = 5326.2856
$ws.Range("N64").Value = -5822.2856

$ws.Range("H67").Value = 5196.8887
$ws.Range("J67").Value = 5326.2856
$ws.Range("L67").Value = 5326.2856
$ws.Range("N67").Value = -7042.2856

$ws.Range("H116").Value = 6659.3335
$ws.Range("I116").Value = 5989.5
$ws.Range("J116").Value = 7999
$ws.Range("K116").Value = 5989.5
$ws.Range("L116").Value = 7999
$ws.Range("M116").Value = -2547.5
$ws.Range("N116").Value = -14883

$ws.Range("H132").Value = 7655.963
$ws.Range("I132").Value = 6977.08
$ws.Range("K132").Value = 20931.24
$ws.Range("M132").Value = -18401.24

$ws.Range("H135").Value = 50584.85
$ws.Range("I135").Value = 499
$ws.Range("K135").Value = 4491
$ws.Range("M135").Value = -1956

$ws.Range("H137").Value = 1949.9615
$ws.Range("I137").Value = 1767.96
$ws.Range("K137").Value = 5303.88
$ws.Range("M137").Value = -2753.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30044.162
$ws.Range("I32").Value = 32285.469
$ws.Range("K32").Value = 32285.469
$ws.Range("M32").Value = -31998.469

$ws.Range("H61").Value = 3996.6
$ws.Range("I61").Value = 3996.6
$ws.Range("K61").Value = 3996.6
$ws.Range("M61").Value = -3784.6

$ws.Range("H97").Value = 647.0833
$ws.Range("I97").Value = 411.6
$ws.Range("K97").Value = 411.6
$ws.Range("M97").Value = 84.39999999999998

$ws.Range("H110").Value = 1349.9286
$ws.Range("I110").Value = 1069.1538
$ws.Range("K110").Value = 1069.1538
$ws.Range("M110").Value = 975.8462

$ws.Range("H136").Value = 3996.6
$ws.Range("I136").Value = 3996.6
$ws.Range("K136").Value = 11989.8
$ws.Range("M136").Value = -9439.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1222
$ws.Range("I94").Value = 1172.3429
$ws.Range("J94").Value = 1415.1111
$ws.Range("K94").Value = 1172.3429
$ws.Range("L94").Value = 1415.1111
$ws.Range("M94").Value = -721.3429000000001
$ws.Range("N94").Value = -2317.1111

$ws.Range("H134").Value = 3515.2917
$ws.Range("I134").Value = 4257.7
$ws.Range("J134").Value = 2985
$ws.Range("K134").Value = 12773.1
$ws.Range("L134").Value = 8955
$ws.Range("M134").Value = -10238.1
$ws.Range("N134").Value = -14025

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 908
$ws.Range("I5").Value = 908
$ws.Range("K5").Value = 908
$ws.Range("M5").Value = -796

$ws.Range("H22").Value = 1323.3334
$ws.Range("I22").Value = 908
$ws.Range("J22").Value = 1738.6666
$ws.Range("K22").Value = 908
$ws.Range("L22").Value = 1738.6666
$ws.Range("M22").Value = -558
$ws.Range("N22").Value = -2438.6666

$ws.Range("H31").Value = 4001.3635
$ws.Range("I31").Value = 2054.2
$ws.Range("K31").Value = 2054.2
$ws.Range("M31").Value = -1759.2

$ws.Range("H34").Value = 4001.3635
$ws.Range("I34").Value = 2054.2
$ws.Range("K34").Value = 2054.2
$ws.Range("M34").Value = -1852.2

$ws.Range("H58").Value = 61054.35
$ws.Range("I58").Value = 101834.5
$ws.Range("J58").Value = 2797
$ws.Range("K58").Value = 101834.5
$ws.Range("L58").Value = 2797
$ws.Range("M58").Value = -101631.5
$ws.Range("N58").Value = -3203

$ws.Range("H64").Value = 69774.5
$ws.Range("I64").Value = 69549
$ws.Range("J64").Value = 70000
$ws.Range("K64").Value = 69549
$ws.Range("L64").Value = 70000
$ws.Range("M64").Value = -69301
$ws.Range("N64").Value = -70496

$ws.Range("H67").Value = 69774.5
$ws.Range("I67").Value = 69549
$ws.Range("J67").Value = 70000
$ws.Range("K67").Value = 69549
$ws.Range("L67").Value = 70000
$ws.Range("M67").Value = -68691
$ws.Range("N67").Value = -71716

$ws.Range("H74").Value = 34397
$ws.Range("I74").Value = 32800
$ws.Range("J74").Value = 34929.332
$ws.Range("K74").Value = 32800
$ws.Range("L74").Value = 34929.332
$ws.Range("M74").Value = -31926
$ws.Range("N74").Value = -36677.332

$ws.Range("H77").Value = 34397
$ws.Range("I77").Value = 32800
$ws.Range("J77").Value = 34929.332
$ws.Range("K77").Value = 98400
$ws.Range("L77").Value = 104787.996
$ws.Range("M77").Value = -94032
$ws.Range("N77").Value = -113523.996

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("L92").ClearContents()

$ws.Range("H96").Value = 2884.4
$ws.Range("J96").Value = 2884.4
$ws.Range("L96").Value = 2884.4
$ws.Range("N96").Value = -8376.4

$ws.Range("H99").Value = 4093.2666
$ws.Range("I99").Value = 3840
$ws.Range("K99").Value = 3840
$ws.Range("M99").Value = -2342

$ws.Range("H105").Value = 1107
$ws.Range("I105").Value = 1107
$ws.Range("K105").Value = 1107
$ws.Range("M105").Value = 640

$ws.Range("H126").Value = 4093.2666
$ws.Range("I126").Value = 3840
$ws.Range("K126").Value = 11520
$ws.Range("M126").Value = -9050

$ws.Range("H136").Value = 61054.35
$ws.Range("I136").Value = 101834.5
$ws.Range("J136").Value = 2797
$ws.Range("K136").Value = 305503.5
$ws.Range("L136").Value = 8391
$ws.Range("M136").Value = -302953.5
$ws.Range("N136").Value = -13491

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.8125
$ws.Range("I2").Value = 20.12
$ws.Range("J2").Value = 87.28570999999999
$ws.Range("K2").Value = 120.72
$ws.Range("L2").Value = 523.71426
$ws.Range("M2").Value = -7.719999999999999
$ws.Range("N2").Value = -749.71426

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H117").Value = 2394.0476
$ws.Range("I117").Value = 1418.6
$ws.Range("J117").Value = 2698.875
$ws.Range("K117").Value = 4255.799999999999
$ws.Range("L117").Value = 8096.625
$ws.Range("M117").Value = -813.7999999999993
$ws.Range("N117").Value = -14980.625

$ws.Range("H131").Value = 11264.091
$ws.Range("J131").Value = 15972.134
$ws.Range("L131").Value = 47916.402
$ws.Range("N131").Value = -57996.402

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 37984.4
$ws.Range("I52").Value = 38307.668
$ws.Range("K52").Value = 38307.668
$ws.Range("M52").Value = -38048.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1598.5714
$ws.Range("I16").Value = 1598.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1598.5714
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -1428.5714
$ws.Range("M16").ClearContents()

$ws.Range("H40").Value = 4498.4346
$ws.Range("I40").Value = 3838.353
$ws.Range("K40").Value = 3838.353
$ws.Range("M40").Value = -3702.353

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("L63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("L66").ClearContents()

$ws.Range("H94").Value = 59550
$ws.Range("I94").Value = 49100
$ws.Range("J94").Value = 70000
$ws.Range("K94").Value = 49100
$ws.Range("L94").Value = 70000
$ws.Range("M94").Value = -48424
$ws.Range("N94").Value = -71352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 171414.33
$ws.Range("I62").Value = 5624.25
$ws.Range("J62").Value = 502994.5
$ws.Range("K62").Value = 5624.25
$ws.Range("L62").Value = 502994.5
$ws.Range("M62").Value = -5000.25
$ws.Range("N62").Value = -504242.5

$ws.Range("H65").Value = 171414.33
$ws.Range("I65").Value = 5624.25
$ws.Range("J65").Value = 502994.5
$ws.Range("K65").Value = 28121.25
$ws.Range("L65").Value = 2514972.5
$ws.Range("M65").Value = -25001.25
$ws.Range("N65").Value = -2521212.5

$ws.Range("H86").Value = 79999
$ws.Range("J86").Value = 79999
$ws.Range("L86").Value = 79999
$ws.Range("N86").Value = -82245

$ws.Range("H89").Value = 79999
$ws.Range("J89").Value = 79999
$ws.Range("L89").Value = 399995
$ws.Range("N89").Value = -411227

$ws.Range("H136").Value = 3356
$ws.Range("I136").Value = 2571.2856
$ws.Range("J136").Value = 6102.5
$ws.Range("K136").Value = 7713.8568
$ws.Range("L136").Value = 18307.5
$ws.Range("M136").Value = -5163.8568
$ws.Range("N136").Value = -23407.5
